# Rename the "Badge No" header in cell A1 to "Badge".
# (Supports an upstream code change that now allows skipping rows when
# reading this data file, so the header column was shortened accordingly.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A1").Value = "Badge"
